$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill F2:F13 with 10, 20, 30, ... 120 and give them the same "red font"
# style already used elsewhere in the sheet (matches cellXfs index 3).
$values = @(10, 20, 30, 40, 50, 60, 70, 80, 90, 100, 110, 120)
for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $cell = $ws.Range("F$row")
    $cell.Value = $values[$i]
    $cell.Font.Color = 255
}

# Update the view: scroll so column B is leftmost and select F15.
$ws.Range("F15").Select()

Write-Host "Applied Wurth World Cup 2026 F2:F13 updates"
